# rekap data 2019.xlsx - "update buku tesis dan tambah data"
#
# 1. Apply the built-in "Comma [0]" cell style (thousands separator, 0
#    decimals -> numFmtId 41) to the three percentage-summary rows
#    (B:N on rows 26, 31 and 36).
# 2. Insert a new (blank) row above row 42, pushing the "agregat" /
#    "positif" / "negatif" / percentage block down by one row
#    (42->43, 43->44, 44->45, 45->46) and growing the used range to
#    A1:W46.
# 3. Leave the final selection on B36:N36, matching the workbook as last
#    saved by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Comma [0] number-format style on the three recap rows ---------
$ws.Range("B26:N26").Style = "Comma [0]"
$ws.Range("B31:N31").Style = "Comma [0]"
$ws.Range("B36:N36").Style = "Comma [0]"

# --- 2. Insert a blank row above row 42 (shifts 42..45 -> 43..46) ------
$ws.Rows(42).Insert() | Out-Null

# --- 3. Restore the selection shown in the saved workbook --------------
$ws.Range("B36:N36").Select() | Out-Null
